$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, 'Total costos y gastos de nómina', 30, 5490109000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(3, 'Aportes al sistema de seguridad social', 31, 659465000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(4, 'Aportes al SENA, ICBF, cajas de compensación', 32, 166548000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(5, 'Efectivo y equivalentes al efectivo', 33, 459434000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(6, 'Inversiones e instrumentos financieros derivados', 34, 9665000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(7, 'Cuentas, documentos y arrendamientos financieros por cobrar', 35, 935790000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(8, 'Inventarios', 36, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(9, 'Activos intangibles', 37, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(10, 'Activos biológicos', 38, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(11, 'Propiedades, planta y equipo, propiedades de inversión y ANCMV', 39, 933615000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(12, 'Otros activos', 40, 149190000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(13, 'Total patrimonio bruto', 41, 2487694000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(14, 'Pasivos', 42, 1365966000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(15, 'Total patrimonio líquido', 43, 1121728000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(16, 'Ingresos brutos de actividades ordinarias', 44, 7337948000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(17, 'Ingresos financieros', 45, 9607000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(18, 'Dividendos y/o participaciones recibidos o capitalizados por sociedades extranjeras (año 2016 y anteriores) o nacionales cualquier año', 46, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(19, 'Dividendos y/o participaciones recibidos por declarantes diferentes a sociedades nacionales, años 2017 y siguientes', 47, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(20, 'Dividendos y/o participaciones recibidos por personas naturales sin residencias fiscal (año 2016 y anteriores)', 48, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(21, 'Otros ingresos', 49, 35950000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(22, 'Total ingresos brutos', 50, 7383505000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(23, 'Devoluciones, rebajas y descuentos en ventas', 51, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(24, 'Ingresos no constitutivos de renta ni ganancia ocasional', 52, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(25, 'Ingresos no constitutivos de dividendos y/o participaciones personas naturales sin residencia fiscal (año 2016 y anteriores)', 53, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(26, 'Total ingresos netos', 54, 7383505000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(27, 'Costos', 55, 6510665000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(28, 'Gastos de administración', 56, 759825000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(29, 'Gastos de distribución y ventas', 57, 40269000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(30, 'Gastos financieros', 58, 30323000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(31, 'Otros gastos y deducciones', 59, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(32, 'Total costos y gastos deducibles', 60, 7341082000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(33, 'Inversiones efectuadas en el año', 61, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(34, 'Inversiones liquidadas de períodos gravables anteriores', 62, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(35, 'Renta Pasiva - ECE sin residencia fiscal en Colombia', 63, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(36, 'Renta líquida ordinaria del ejercicio sin casilla 47 y 48', 64, 42423000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(37, 'Pérdida líquida del ejercicio sin casilla 47 y 48', 65, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(38, 'Compensaciones', 66, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(39, 'Renta líquida sin casilla 47 y 48', 67, 42423000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(40, 'Renta presuntiva', 68, 42749000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(41, 'Renta exenta', 69, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(42, 'Rentas gravables', 70, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(43, 'Sin dividendos gravados al 5%, 35% y 33% de personas naturales sin residencia fisca', 71, 42749000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(44, 'Dividendos gravados a la tarifa del 5%', 72, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(45, 'Dividendos gravados a la tarifa del 35%', 73, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(46, 'Dividendos gravados, a la tarifa del 33% Personas Naturales sin residencia fiscal (año 2016 y anteriores)', 74, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(47, 'Ingresos por ganancias ocasionales', 75, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(48, 'Costos por ganancias ocasionales', 76, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(49, 'Ganancias ocasionales no gravadas y exentas', 77, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(50, 'Ganancias ocasionales gravables', 78, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(51, 'Impuesto sobre la renta líquida gravable', 79, 14107000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(52, 'Descuentos tributarios', 80, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(53, 'Impuesto neto de renta', 81, 14107000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(54, 'Sobretasa', 82, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(55, 'Impuesto de ganancias ocasionales', 83, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(56, 'Descuento por impuestos pagados en el exterior por ganancías ocasionales', 84, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(57, 'Impuesto dividendos gravados a la tarifa del 5%', 85, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(58, 'Impuesto dividendos gravados a la tarifa del 35%', 86, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(59, 'Impuesto dividendos gravados a la tarifa del 33%', 87, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(60, 'Total impuesto a cargo', 88, 14107000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(61, 'Valor inversion obras por impuestos hasta del 50% del valor de la casilla 88 (Modalidad de pago 1)', 89, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(62, 'Descuento efectivo inversión obras por impuestos (Modalidad de pago 2)', 90, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(63, 'Anticipo renta liquidado año gravable anterior', 91, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(64, 'Anticipo sobretasa liquidado año gravable anterior', 92, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(65, 'Saldo a favor año gravable anterior sin solicitud de devolución y/o compensación', 93, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(66, 'Autorretenciones', 94, 59176000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(67, 'Otras retenciones', 95, 21289000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(68, 'Total retenciones año gravable a declarar', 96, 80465000, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(69, 'Anticipo renta para el año gravable siguiente', 97, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(70, 'Saldo a pagar por impuesto', 98, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(71, 'Sanciones', 99, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(72, 'Total saldo a pagar', 100, 0, 2018, 'SEGURIDAD IMPERIO LTDA.'),
    @(73, 'Total saldo a favor', 101, 66358000, 2018, 'SEGURIDAD IMPERIO LTDA.')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
